$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: bump "Last Updated" timestamp
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 07:04 PM"

# ------------------------------------------------------------------
# 2. "Top Losers" sheet: update a few Weekly (column D) values
# ------------------------------------------------------------------
$losers = $wb.Worksheets.Item("Top Losers")
$losers.Range("D18").Value = 5.978
$losers.Range("D48").Value = -2.9654
$losers.Range("D56").Value = 5.2953

# ------------------------------------------------------------------
# 3. Add new "distance from Dma50" sheet at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "distance from Dma50"

# Header row
$ws.Cells.Item(1, 1).Value = "Icon"
$ws.Cells.Item(1, 2).Value = "Stock"
$ws.Cells.Item(1, 3).Value = "Distance From Sma50"

# Copy the header formatting (bold, bordered, centered) from an existing sheet
# so the new header matches the workbook's other tables, without disturbing
# any other sheet content.
$wb.Worksheets.Item("Top Losers").Range("A1:C1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

# Data rows (Icon, Stock, Distance From Sma50) sorted high -> low, as in source
$data = @(
    @("📈", "NIFTYPSUBANK", 10.2033),
    @("📈", "NIFTYMETAL", 8.624700000000001),
    @("📈", "NIFTYOILANDGAS", 6.396),
    @("📈", "NIFTYCOMMODITIES", 5.7207),
    @("📈", "CNXINFRA", 5.6012),
    @("📈", "CNXREALTY", 5.4493),
    @("📈", "NIFTYPVTBANK", 5.0059),
    @("📈", "BANKNIFTY", 4.9192),
    @("📈", "NIFTYFINSERVICE", 3.9783),
    @("📈", "NIFTYMIDCAP50", 3.9228),
    @("📈", "NIFTY", 3.7191),
    @("📈", "CNXENERGY", 3.706),
    @("📈", "CNXMIDCAP", 3.6313),
    @("📈", "NIFTY200", 3.5915),
    @("📈", "NIFTY100", 3.5759),
    @("📈", "NIFTY500", 3.3379),
    @("📈", "CNXSMALLCAP", 2.9205),
    @("📈", "NIFTY50VALUE20", 2.893),
    @("📈", "NIFTYCPSE", 2.837),
    @("📈", "CNXNIFTYJUNIOR", 2.8291),
    @("📈", "NIFTYHEALTHCARE", 2.162),
    @("📈", "CNXIT", 2.0641),
    @("📈", "NIFTYCONSUMPTION", 2.0245),
    @("📈", "CNXPHARMA", 1.573),
    @("📈", "NIFTYAUTO", 1.5538),
    @("📈", "NIFTYGROWSECT15", 1.5325),
    @("📈", "NIFTYFMCG", 1.3194),
    @("📈", "NIFTYCONSURDURBL", 0.4031),
    @("📈", "NIFTYMEDIA", -1.9217)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Match page margins used throughout the rest of the workbook (in points:
# 0.75in/0.75in/1in/1in/0.5in/0.5in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Leave the first sheet focused, matching the workbook's original state.
$wb.Worksheets.Item("Metadata").Activate()
